# Movie_Data.xlsx edit: rename director sheet, reposition it, and populate
# it with the per-director "highest gross film" rollup.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # MovieInfo
$ws2 = $wb.Worksheets.Item(2)   # "Director's Highest gross Films" -> Directors

# Recreate the second sheet so it gets a fresh sheetId (matches a sheet that
# was deleted and re-added rather than merely renamed), keeping it right
# after MovieInfo.
$ws2.Delete()
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Directors"

$data = @(
    @('Antoine Fuqua', 'King Arthur', 203567857),
    @('Michael Chaves', $null, 0),
    @('Cal Brunker', 'PAW Patrol: The Movie', 144327371),
    @('Kenneth Branagh', 'Thor: The Dark World', 644783140),
    @('Kevin Greutert', 'Titanic', 2264743305),
    @('Nia Vardalos', 'My Big Fat Greek Wedding', 368744044),
    @('Gareth Edwards', 'Rogue One: A Star Wars Story', 1058682142),
    @('Scott Waugh', 'Need for Speed', 203277636),
    @('Atlee', 'Jawan', 27514156),
    @('Andrew Hyatt', 'Paul, Apostle of Christ', 25915966),
    @('Craig Gillespie', 'Cruella', 233503234),
    @('Bishal Dutta', 'It Lives Inside', 5482605),
    @('Jonathan Demme', 'Philadelphia', 206678440),
    @('Matthew Crouch', 'Route 60: The Biblical Highway', 1606339),
    @('Wuershan', 'The Butcher, the Chef, and the Swordsman', 2089266),
    @('Sean Olson', 'A Question of Faith', 2587072),
    @('Miles Joris-Peyrafitte', 'The Good Mother', 503378),
    @('Rudy Valdez', 'Carlos', 423024),
    @('Aitch Alberto', 'Aristotle and Dante Discover the Secrets of the Universe', 407838),
    @('Michael A. Goorjian', 'Amerikatsi', 376719),
    @('Michael Jai White', 'Outlaw Johnny Black', 319848),
    @('Jim Capobianco', 'The Lion King', 968511805),
    @('Pierre-Luc Granjon', 'Zibilla', 88667),
    @('Tyler Sansom', $null, 0),
    @('Kaige Chen', 'The Battle at Lake Changjin', 902548476),
    @('Nicol Paone', 'The Kill Room', 617143),
    @('Peter Facinelli', 'On Fire', 205118),
    @('Nick Lyon', 'On Fire', 205118),
    @('Peter Lepeniotis', 'Toy Story 2', 497375381),
    @('Joshua Tickell', 'Common Ground', 88910),
    @('Rebecca Harrell Tickell', 'Pump', 89787),
    @('Sébastien Marnier', 'The Origin of Evil', 1120899),
    @('Bethann Hardison', $null, 0),
    @('Frédéric Tcheng', 'Shortbus', 5557564),
    @('Klaus Härö', 'Letters to Father Jacob', 1332577),
    @('Stephen Gyllenhaal', 'Uncharitable', 38683),
    @('Danny O''Malley', 'Amanda Knox', 397),
    @('Alex Rivest', $null, 0),
    @('Richard Dewey', 'Radical Wolfe', 32535),
    @('Aristotle Torres', 'Story Ave', 25570),
    @('Hiroshi Akabane', $null, 0),
    @('Luca Balser', 'Uncut Gems', 50023780),
    @('Adil El Arbi', 'Black', 1692776),
    @('Bilall Fallah', 'Bad Boys for Life', 426505244),
    @('John Stalberg Jr.', 'High School', 221590),
    @('Jude Okwudiafor Johnson', 'Senior Year: Love Never Fails', 5300),
    @('Cru Ennis', 'Boys of Abu Ghraib', 62096),
    @('Lee Roy Kunz', 'Deliver Us', 4883),
    @('Andrea Di Stefano', 'Escobar: Paradise Lost', 6760531),
    @('Miguel Ángel Vivas', 'Extinction', 2350695),
    @('Frank Cimière', $null, 0)
)

$ws2.Cells.Item(1,1).Value = "Film Director"
$ws2.Cells.Item(1,2).Value = "Film Title"
$ws2.Cells.Item(1,3).Value = "Worldwide Gross"

$r = 2
foreach ($row in $data) {
    $ws2.Cells.Item($r,1).Value = $row[0]
    if ($row[1] -ne $null) {
        $ws2.Cells.Item($r,2).Value = $row[1]
    }
    $ws2.Cells.Item($r,3).Value = $row[2]
    $r++
}

$ws2.Columns.Item(3).ColumnWidth = 23.5703125

# Selection / active-sheet state: MovieInfo ends with a single-cell
# selection on E1, Directors ends up the active (front) tab selected at G15.
$ws1.Range("E1").Select()
$ws2.Activate()
$ws2.Range("G15").Select()
